$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the full contents of row 15 and row 16 (every column that
# differs between the two rows, including a handful of cells that exist only
# on row 15 and must move over to row 16).

# 1) Columns present (and differing) on both rows 15 and 16 -> straight swap.
$swapCols = @("A","B","E","F","G","H","Q","R")
foreach ($col in $swapCols) {
    $c15 = $ws.Range($col + "15")
    $c16 = $ws.Range($col + "16")
    $v15 = $c15.Value2
    $v16 = $c16.Value2
    $c15.Value2 = $v16
    $c16.Value2 = $v15
}

# 2) Columns that only exist on row 15 in the "before" state and must move to
#    row 16 ("after" state) - J, K, L, M, N, AF, AO.
#    J, K, L, N, AF carry no text (they are present-but-blank cells on row 15);
#    M and AO carry real text. Touching a destination cell's formatting with a
#    no-op assignment is enough to materialize the (still blank) cell so the
#    row-16 shape mirrors what row-15 used to look like.
$blankMoveCols = @("J","K","L","N","AF")
foreach ($col in $blankMoveCols) {
    $dst = $ws.Range($col + "16")
    $dst.Font.Bold = $false
    $ws.Range($col + "15").ClearContents()
}

$textMoveCols = @("M","AO")
foreach ($col in $textMoveCols) {
    $src = $ws.Range($col + "15")
    $dst = $ws.Range($col + "16")
    $dst.Value2 = $src.Value2
    $src.ClearContents()
}
